$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow Table1 by one row (A1:I7 -> A1:I8)
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Carry the existing per-column formatting down into the new row 8
# (so it gets the same banded/highlighted styling as the rows above)
$ws.Range("C7:F7").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

# New "Geothermal" generator row
$ws.Range("A8").Value2 = "Geothermal"
$ws.Range("B8").Value2 = 38816715
$ws.Range("C8").Value2 = 0.3
$ws.Range("D8").Formula = "=(1-Table1[[#This Row],[Direct Pay ITC]])*Table1[[#This Row],[$/MW Gross Capital Cost]]"
$ws.Range("E8").Value2 = 0.029
$ws.Range("F8").Value2 = 30
$ws.Range("G8").Formula = "=PMT(Table1[[#This Row],[Annual Rate]],Table1[[#This Row],[Term]],Table1[[#This Row],[Net Capital Cost]])"

# Helper calculations added below the table
$ws.Range("B12").Formula = "=B8/1000"
$ws.Range("G15").Formula = "=495631/2"
$ws.Range("G17").Formula = "=G8/(0.8*8760)"
$ws.Range("G20").Formula = "=G15-G17"

# B8 uses a plain currency format (distinct from the accounting style used
# by the other $/MW Gross Capital Cost cells above it)
$ws.Range("B8").NumberFormat = '"$"#,##0.00'

# Widen column G to fit the new values, and drop the old "best fit" width
$ws.Columns("G").ColumnWidth = 15.1667

# Match the saved selection
$ws.Range("G15").Select() | Out-Null
